$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the Date value (row 8, column B) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-01T18:52:33+00:00"

# --- Concepts sheet: rename the TEMPORARY concept's display text ---
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("C2").Value = "Temporary absence"

# --- Concepts sheet: add a new NEVER / "Permanent absence" concept row ---
# Set the values first (leading "'" keeps "1" stored as text, like A2).
$concepts.Range("A3").Value = "'1"
$concepts.Range("B3").Value = "NEVER"
$concepts.Range("C3").Value = "Permanent absence"

# Copy the formatting from row 2 onto row 3 after the values are in place,
# so the new row's style matches (and doesn't pick up a stray quote-prefix style).
$concepts.Range("A2:D2").Copy()
$concepts.Range("A3:D3").PasteSpecial(-4122)
$excel.CutCopyMode = 0
